# Apply Danish localization changes to the report workbook.
# 1) Rename the two worksheets from English to Danish.
# 2) Replace the English text values with their Danish equivalents
#    wherever they occur (headers, farm name, group labels, etc).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("19 - 03. Control floating laye")
$ws2 = $wb.Worksheets.Item("23 - 04. Feeding documentation")

# --- Rename worksheets ---
$ws1.Name = "19 - 03. Kontrol flydelag"
$ws2.Name = "23 - 04. Foderindlægssedler"

# --- Sheet 1: "19 - 03. Kontrol flydelag" ---

# Header row
$ws1.Range("B1").Value = "Ejendom"
$ws1.Range("C1").Value = "Dato"
$ws1.Range("D1").Value = "Udført af"
$ws1.Range("E1").Value = "Område"
$ws1.Range("F1").Value = "Flydelag OK"
$ws1.Range("G1").Value = "Vælg årsag til manglende flydelag"
$ws1.Range("H1").Value = "Kommentar"

# Farm name "Tjørntved" -> "Farm 1"
$ws1.Range("B2").Value = "Farm 1"
$ws1.Range("B3").Value = "Farm 1"

# Group labels "G1/G2/G3: Floating layer" -> "G1/G2/G3: Flydelag"
$ws1.Range("E2").Value = "G1: Flydelag"
$ws1.Range("E3").Value = "G2: Flydelag"
$ws1.Range("E4").Value = "G1: Flydelag"
$ws1.Range("E5").Value = "G2: Flydelag"
$ws1.Range("E6").Value = "G3: Flydelag"

# "Slurry tank empty" -> "Beholder tom"
$ws1.Range("G4").Value = "Beholder tom"
$ws1.Range("G5").Value = "Beholder tom"

# --- Sheet 2: "23 - 04. Foderindlægssedler" ---

# Header row
$ws2.Range("B1").Value = "Ejendom"
$ws2.Range("C1").Value = "Dato"
$ws2.Range("D1").Value = "Udført af"
$ws2.Range("E1").Value = "Område"
$ws2.Range("F1").Value = "Kommentar"

# Farm name "Tjørntved" -> "Farm 1"
$ws2.Range("B2").Value = "Farm 1"
